$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number: force Text format so
# Excel COM does not silently convert the literal string into a numeric
# value (which would also lose the original "sig-figs" text formatting).
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D8",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D15",
    "D18",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D39",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D51"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "65.916.72"
$ws.Range("D3").Value = "3.278.31"
$ws.Range("E3").Value = "  -4.90%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "557.28"
$ws.Range("E5").Value = "  -3.16%  "
$ws.Range("D6").Value = "185.05"
$ws.Range("E6").Value = "  -2.38%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("D9").Value = "3.275.03"
$ws.Range("E9").Value = "  -4.51%  "
$ws.Range("D10").Value = "0.187"
$ws.Range("E10").Value = "  -7.82%  "
$ws.Range("D11").Value = "0.587"
$ws.Range("E11").Value = "  -4.36%  "
$ws.Range("D12").Value = "47.37"
$ws.Range("E12").Value = "  -7.12%  "
$ws.Range("D13").Value = "0.0000267"
$ws.Range("E13").Value = "  -5.74%  "
$ws.Range("D14").Value = "8.59"
$ws.Range("E14").Value = "  -4.96%  "
$ws.Range("D15").Value = "631.53"
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "3.800.72"
$ws.Range("E16").Value = "  -4.54%  "
$ws.Range("D17").Value = "65.916.03"
$ws.Range("E17").Value = "  -4.07%  "
$ws.Range("D18").Value = "17.90"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("D20").Value = "3.265.31"
$ws.Range("E20").Value = "  -4.89%  "
$ws.Range("D21").Value = "11.37"
$ws.Range("E21").Value = "  -6.54%  "
$ws.Range("D22").Value = "0.905"
$ws.Range("E22").Value = "  -3.39%  "
$ws.Range("D23").Value = "17.80"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "107.06"
$ws.Range("E24").Value = "  +9.05%  "
$ws.Range("D25").Value = "4.95"
$ws.Range("E25").Value = "  -6.92%  "
$ws.Range("D26").Value = "3.98"
$ws.Range("D27").Value = "2.67"
$ws.Range("E27").Value = "  -5.69%  "
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("D29").Value = "8.70"
$ws.Range("E29").Value = "  -5.13%  "
$ws.Range("D30").Value = "30.44"
$ws.Range("E30").Value = "  -5.27%  "
$ws.Range("D31").Value = "4.05"
$ws.Range("E31").Value = "  -5.76%  "
$ws.Range("D32").Value = "6.29"
$ws.Range("E32").Value = "  -5.50%  "
$ws.Range("D33").Value = "11.03"
$ws.Range("E33").Value = "  -4.06%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.105"
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "541.99"
$ws.Range("E35").Value = "  +11.22%  "
$ws.Range("E36").Value = "  -5.90%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "3.669.80"
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("D39").Value = "3.43"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "0.0₃0731"
$ws.Range("E40").Value = "  -6.70%  "
$ws.Range("D41").Value = "0.131"
$ws.Range("E41").Value = "  -0.76%  "
$ws.Range("D42").Value = "2.73"
$ws.Range("E42").Value = "  -5.54%  "
$ws.Range("D43").Value = "3.33"
$ws.Range("E43").Value = "  -5.27%  "
$ws.Range("D44").Value = "32.62"
$ws.Range("E44").Value = "  -4.14%  "
$ws.Range("D45").Value = "0.338"
$ws.Range("E45").Value = "  -7.94%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0414"
$ws.Range("E46").Value = "  -4.70%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "3.25"
$ws.Range("E47").Value = "  -2.23%  "
$ws.Range("D48").Value = "2.62"
$ws.Range("E48").Value = "  -5.66%  "
$ws.Range("E49").Value = "  -2.99%  "
$ws.Range("D51").Value = "1.25"
$ws.Range("E51").Value = "  +1.95%  "

# Drop the temporary Text number-format again so the cells end up with
# no explicit style, matching the original (unstyled) inline-string cells.
foreach ($ref in $textCells) {
    $ws.Range($ref).ClearFormats()
}
